# fix(source): solve numeric formatting problem
#
# The "Soft Min" / "Soft Max" columns (I/J) on the Plants sheet held a
# handful of values that had been typed with a comma decimal separator
# (e.g. "0,2"), so Excel stored them as plain text instead of numbers.
# Convert those text cells to real numeric values. Also correct a few
# "Tune" flags (column E) that were left as "no" and should be "yes", and
# highlight the two still-empty "Hard Max" cells (K2, K3) in yellow so
# they are easy to spot for follow-up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plants")

# --- Cryptomonad block -----------------------------------------------
$ws.Range("I2").Value = 0.2
$ws.Range("I3").Value = 0.04
$ws.Range("J3").Value = 0.1
$ws.Range("I4").Value = 0.001
$ws.Range("J4").Value = 1.52
$ws.Range("I9").Value = 0.005
$ws.Range("J9").Value = 0.01
$ws.Range("I10").Value = 0.07
$ws.Range("J10").Value = 0.32
$ws.Range("I12").Value = 0.05
$ws.Range("J12").Value = 0.2

# --- Cyclotella Nana block --------------------------------------------
$ws.Range("I15").Value = 0.34
$ws.Range("J15").Value = 3.4
$ws.Range("I16").Value = 0.001
$ws.Range("J16").Value = 1.55

# --- Fix "Tune" flags that were mistakenly left as "no" ---------------
$ws.Range("E6").Value = "yes"
$ws.Range("E7").Value = "yes"
$ws.Range("E11").Value = "yes"
$ws.Range("E19").Value = "yes"
$ws.Range("E20").Value = "yes"
$ws.Range("E21").Value = "yes"

# --- Highlight the still-blank Hard Max cells in yellow ----------------
$ws.Range("K2:K3").Interior.Color = 65535

# --- Restore the view: scroll back to the top, select the Tune column --
[void]$ws.Range("E5:E21").Select()
